$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting of the last existing row (31) into the new row (32)
$ws.Range("A31:D31").Copy()
$ws.Range("A32:D32").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Rows("32:32").RowHeight = $ws.Rows("31:31").RowHeight

# Fill in the new July 1st, 2020 data row
$ws.Range("A32").Value = 30
$ws.Range("B32").Value = "2020-07-01"
$ws.Range("C32").Value = 1
$ws.Range("D32").Value = 4.7

# Update the active selection to match the new last row
[void]$ws.Activate()
[void]$ws.Range("B33").Select()
